# Update the "想去人数" (Want-to-go count) column F for several rows
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All
# Types) sheets, matching the refreshed scrape data from commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- 展览 sheet ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 35
$ws.Range("F6").Value = 8897
$ws.Range("F8").Value = 238
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 400
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 1905
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 36
$ws.Range("F31").Value = 81
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 1031
$ws.Range("F35").Value = 15
$ws.Range("F36").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 1194
$ws.Range("F41").Value = 649
$ws.Range("F42").Value = 0
$ws.Range("F44").Value = 1100
$ws.Range("F46").Value = 976
$ws.Range("F47").Value = 1358
$ws.Range("F48").Value = 0

# ---- 演出 sheet ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 37
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = 3
$ws.Range("F9").Value = 42
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F15").Value = 94
$ws.Range("F17").Value = 7
$ws.Range("F19").Value = 0

# ---- 全部类型 sheet ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 37
$ws.Range("F3").Value = 35
$ws.Range("F6").Value = 15
$ws.Range("F7").Value = 1145
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 7107
$ws.Range("F13").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 1094
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 272
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 104
$ws.Range("F26").Value = 190
$ws.Range("F28").Value = 83
$ws.Range("F30").Value = 0
$ws.Range("F33").Value = 81
$ws.Range("F35").Value = 1031
$ws.Range("F38").Value = 2047
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 650
$ws.Range("F43").Value = 106
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("F47").Value = 976
